$wb = $excel.ActiveWorkbook

# --- Update the "Conversión del día" note on Hoja1!A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$cellA1 = $wsHoja1.Range("A1")
$oldText = $cellA1.Value()
$newText = $oldText.Replace("1000 Bs = 3.68 = 14121.32 pesos", "1000 Bs = 3.65 = 13949.41 pesos").Replace("14121.32 pesos = 3.66 = 932.68 Bs", "13949.41 pesos = 3.62 = 944.23 Bs")
$cellA1.Value = $newText

# --- Update rate figures on "tasas" sheet ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 273.99
$wsTasas.Range("O10").Value = 3822
$wsTasas.Range("N12").Value = 3850
$wsTasas.Range("O12").Value = 260.605
